$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("W2").Value = 2.52

# Row 3 updates
$ws.Range("J3").Value = 4.4
$ws.Range("Q3").Value = 1.59
$ws.Range("T3").Value = 1.85
$ws.Range("U3").Value = 1.96
$ws.Range("AF3").Value = 95
